$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ecm1"
$ws.Cells.Item(2, 3).Value = "Itgb4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 11.001142
$ws.Cells.Item(2, 8).Value = 33.003426
$ws.Cells.Item(2, 9).Value = 0.04976752323647229
$ws.Cells.Item(2, 10).Value = 0.04976752323647229
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 4.187598666666667
$ws.Cells.Item(2, 14).Value = 12.562796
$ws.Cells.Item(2, 15).Value = 0.4284941686600627
$ws.Cells.Item(2, 16).Value = 0.4284941686600626
$ws.Cells.Item(2, 17).Value = 46.06836757101067
$ws.Cells.Item(2, 18).Value = 414.615308139096
$ws.Cells.Item(2, 19).Value = 0.02132509349548255
$ws.Cells.Item(2, 20).Value = 0.02132509349548254

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ecm1"
$ws.Cells.Item(3, 3).Value = "Itgb4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 11.001142
$ws.Cells.Item(3, 8).Value = 33.003426
$ws.Cells.Item(3, 9).Value = 0.04976752323647229
$ws.Cells.Item(3, 10).Value = 0.04976752323647229
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.204585666666667
$ws.Cells.Item(3, 14).Value = 3.613757
$ws.Cells.Item(3, 15).Value = 0.1232586918910792
$ws.Cells.Item(3, 16).Value = 0.1232586918910792
$ws.Cells.Item(3, 17).Value = 13.25181797016467
$ws.Cells.Item(3, 18).Value = 119.266361731482
$ws.Cells.Item(3, 19).Value = 0.006134279812786463
$ws.Cells.Item(3, 20).Value = 0.006134279812786462

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ecm1"
$ws.Cells.Item(4, 3).Value = "Itgb4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 11.001142
$ws.Cells.Item(4, 8).Value = 33.003426
$ws.Cells.Item(4, 9).Value = 0.04976752323647229
$ws.Cells.Item(4, 10).Value = 0.04976752323647229
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.380641
$ws.Cells.Item(4, 14).Value = 13.141923
$ws.Cells.Item(4, 15).Value = 0.4482471394488581
$ws.Cells.Item(4, 16).Value = 0.4482471394488581
$ws.Cells.Item(4, 17).Value = 48.19205369202199
$ws.Cells.Item(4, 18).Value = 433.728483228198
$ws.Cells.Item(4, 19).Value = 0.02230814992820328
$ws.Cells.Item(4, 20).Value = 0.02230814992820328

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ecm1"
$ws.Cells.Item(5, 3).Value = "Itgb4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 200.8411913333333
$ws.Cells.Item(5, 8).Value = 602.5235739999999
$ws.Cells.Item(5, 9).Value = 0.9085755512039061
$ws.Cells.Item(5, 10).Value = 0.9085755512039061
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 4.187598666666667
$ws.Cells.Item(5, 14).Value = 12.562796
$ws.Cells.Item(5, 15).Value = 0.4284941686600627
$ws.Cells.Item(5, 16).Value = 0.4284941686600626
$ws.Cells.Item(5, 17).Value = 841.0423050392116
$ws.Cells.Item(5, 18).Value = 7569.380745352903
$ws.Cells.Item(5, 19).Value = 0.389319325477976
$ws.Cells.Item(5, 20).Value = 0.3893193254779759

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ecm1"
$ws.Cells.Item(6, 3).Value = "Itgb4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 200.8411913333333
$ws.Cells.Item(6, 8).Value = 602.5235739999999
$ws.Cells.Item(6, 9).Value = 0.9085755512039061
$ws.Cells.Item(6, 10).Value = 0.9085755512039061
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.204585666666667
$ws.Cells.Item(6, 14).Value = 3.613757
$ws.Cells.Item(6, 15).Value = 0.1232586918910792
$ws.Cells.Item(6, 16).Value = 0.1232586918910792
$ws.Cells.Item(6, 17).Value = 241.9304203563909
$ws.Cells.Item(6, 18).Value = 2177.373783207518
$ws.Cells.Item(6, 19).Value = 0.1119898339256097
$ws.Cells.Item(6, 20).Value = 0.1119898339256097

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ecm1"
$ws.Cells.Item(7, 3).Value = "Itgb4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 200.8411913333333
$ws.Cells.Item(7, 8).Value = 602.5235739999999
$ws.Cells.Item(7, 9).Value = 0.9085755512039061
$ws.Cells.Item(7, 10).Value = 0.9085755512039061
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.380641
$ws.Cells.Item(7, 14).Value = 13.141923
$ws.Cells.Item(7, 15).Value = 0.4482471394488581
$ws.Cells.Item(7, 16).Value = 0.4482471394488581
$ws.Cells.Item(7, 17).Value = 879.8131572436446
$ws.Cells.Item(7, 18).Value = 7918.318415192802
$ws.Cells.Item(7, 19).Value = 0.4072663918003204
$ws.Cells.Item(7, 20).Value = 0.4072663918003204

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ecm1"
$ws.Cells.Item(8, 3).Value = "Itgb4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.208289333333333
$ws.Cells.Item(8, 8).Value = 27.624868
$ws.Cells.Item(8, 9).Value = 0.04165692555962159
$ws.Cells.Item(8, 10).Value = 0.04165692555962159
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 4.187598666666667
$ws.Cells.Item(8, 14).Value = 12.562796
$ws.Cells.Item(8, 15).Value = 0.4284941686600627
$ws.Cells.Item(8, 16).Value = 0.4284941686600626
$ws.Cells.Item(8, 17).Value = 38.56062013454756
$ws.Cells.Item(8, 18).Value = 347.045581210928
$ws.Cells.Item(8, 19).Value = 0.01784974968660417
$ws.Cells.Item(8, 20).Value = 0.01784974968660417

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ecm1"
$ws.Cells.Item(9, 3).Value = "Itgb4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.208289333333333
$ws.Cells.Item(9, 8).Value = 27.624868
$ws.Cells.Item(9, 9).Value = 0.04165692555962159
$ws.Cells.Item(9, 10).Value = 0.04165692555962159
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.204585666666667
$ws.Cells.Item(9, 14).Value = 3.613757
$ws.Cells.Item(9, 15).Value = 0.1232586918910792
$ws.Cells.Item(9, 16).Value = 0.1232586918910792
$ws.Cells.Item(9, 17).Value = 11.09217334545289
$ws.Cells.Item(9, 18).Value = 99.829560109076
$ws.Cells.Item(9, 19).Value = 0.00513457815268302
$ws.Cells.Item(9, 20).Value = 0.005134578152683019

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ecm1"
$ws.Cells.Item(10, 3).Value = "Itgb4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.208289333333333
$ws.Cells.Item(10, 8).Value = 27.624868
$ws.Cells.Item(10, 9).Value = 0.04165692555962159
$ws.Cells.Item(10, 10).Value = 0.04165692555962159
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.380641
$ws.Cells.Item(10, 14).Value = 13.141923
$ws.Cells.Item(10, 15).Value = 0.4482471394488581
$ws.Cells.Item(10, 16).Value = 0.4482471394488581
$ws.Cells.Item(10, 17).Value = 40.33820979346267
$ws.Cells.Item(10, 18).Value = 363.043888141164
$ws.Cells.Item(10, 19).Value = 0.0186725977203344
$ws.Cells.Item(10, 20).Value = 0.0186725977203344
